# Update Excel file with latest predictions.
#
# The site's odds feed moved forward by one day: matches that were
# previously dated 31-12-2024 have finished/expired and are removed from
# each market sheet, with the remaining (already 01-01-2025) rows shifting
# up to fill the gap.

$wb = $excel.ActiveWorkbook

# "Draw" sheet: drop the single outdated 31-12-2024 match (row 2).
$wsDraw = $wb.Worksheets.Item("Draw")
$wsDraw.Rows.Item(2).Delete()

# "Btts" sheet: drop the single outdated 31-12-2024 match (row 2).
$wsBtts = $wb.Worksheets.Item("Btts")
$wsBtts.Rows.Item(2).Delete()

# "Over_Under" sheet: drop the two outdated 31-12-2024 matches (rows 2 & 3).
$wsOverUnder = $wb.Worksheets.Item("Over_Under")
$wsOverUnder.Rows.Item(2).Delete()
$wsOverUnder.Rows.Item(2).Delete()
